# "Da gestão de negócios ao mérito Científico"
#   -> "Da gestão e inteligência de negócios ao mérito Científico"
#
# The run immediately after the edited text ("o mérito Científico") must
# stay untouched, so before mutating any text we temporarily mark it with
# a distinct character property (Bold) - this stops the engine from
# coalescing it into the run we are about to edit. We then insert the new
# words, remove the temporary marker, and finally re-split the edited
# run into the four pieces by toggling Bold on/off (which cleanly leaves
# no formatting trace) across each boundary.

$d = $word.ActiveDocument

$full = $d.Content
$full.Find.Execute("Da gestão de negócios a", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$startPos = $full.Start
$endPos   = $full.End

# 1) Mark the following run ("o mérito Científico", 19 chars) so it is
#    not pulled into the edited run.
$tail = $d.Range($endPos, $endPos + 19)
$tail.Bold = 1

# 2) Insert " e inteligência" right after "Da gestão".
$splitPos = $startPos + 9
$insertionPoint = $d.Range($splitPos, $splitPos)
$insertionPoint.InsertAfter(" e inteligência")

# 3) Un-mark the tail run (its offsets shifted by the 15 inserted chars).
$newTailStart = $endPos + 15
$tail2 = $d.Range($newTailStart, $newTailStart + 19)
$tail2.Bold = 0

# 4) Re-split the edited run into "Da gestão" | " e i" | "nteligência" |
#    " de negócios a" by toggling Bold on/off across each boundary.
$b1 = $startPos
$b2 = $startPos + 9
$b3 = $startPos + 13
$b4 = $startPos + 24
$b5 = $startPos + 38

$p1 = $d.Range($b1, $b2); $p1.Bold = 1; $p1.Bold = 0
$p2 = $d.Range($b2, $b3); $p2.Bold = 1; $p2.Bold = 0
$p3 = $d.Range($b3, $b4); $p3.Bold = 1; $p3.Bold = 0
$p4 = $d.Range($b4, $b5); $p4.Bold = 1; $p4.Bold = 0
